$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

function Set-TextValue([string]$cellRef, [string]$text) {
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextValue "D2" "67.951.05"
$ws.Range("E2").Value = "  -1.85%  "
Set-TextValue "D3" "3.268.04"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "185.89"
$ws.Range("E5").Value = "  -0.48%  "
Set-TextValue "D6" "580.91"
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue "D8" "0.600"
$ws.Range("E8").Value = "  -1.01%  "
Set-TextValue "D9" "3.269.02"
$ws.Range("E9").Value = "  -1.49%  "
Set-TextValue "D10" "0.130"
$ws.Range("E10").Value = "  -5.41%  "
Set-TextValue "D11" "6.54"
$ws.Range("E11").Value = "  -2.77%  "
Set-TextValue "D12" "0.410"
$ws.Range("E12").Value = "  -3.32%  "
Set-TextValue "D13" "3.828.02"
$ws.Range("E13").Value = "  -1.75%  "
Set-TextValue "D14" "0.137"
$ws.Range("E14").Value = "  +0.10%  "
Set-TextValue "D15" "27.49"
$ws.Range("E15").Value = "  -6.48%  "
Set-TextValue "D16" "67.998.29"
$ws.Range("E16").Value = "  -1.82%  "
Set-TextValue "D17" "0.0000168"
$ws.Range("E17").Value = "  -3.79%  "
Set-TextValue "D18" "3.308.05"
$ws.Range("E18").Value = "  -0.21%  "
Set-TextValue "D19" "5.73"
$ws.Range("E19").Value = "  -3.36%  "
Set-TextValue "D20" "13.54"
$ws.Range("E20").Value = "  -2.20%  "
Set-TextValue "D21" "398.79"
$ws.Range("E21").Value = "  +1.17%  "
Set-TextValue "D22" "7.61"
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("E23").Value = "  +0.19%  "
Set-TextValue "D24" "71.34"
$ws.Range("E24").Value = "  -1.18%  "
Set-TextValue "D25" "0.511"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("E26").Value = "  -5.01%  "
$ws.Range("E27").Value = "  -0.62%  "
Set-TextValue "D28" "9.51"
$ws.Range("E28").Value = "  -3.74%  "
Set-TextValue "D29" "1.01"
$ws.Range("E29").Value = "  +0.98%  "
Set-TextValue "D30" "1.95"
$ws.Range("E30").Value = "  -2.94%  "
Set-TextValue "D31" "22.68"
$ws.Range("E31").Value = "  -2.24%  "
Set-TextValue "D32" "5.50"
$ws.Range("E32").Value = "  -7.13%  "
Set-TextValue "D33" "6.95"
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("E34").Value = "  -5.92%  "
Set-TextValue "D36" "162.79"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("E37").Value = "  -6.36%  "
Set-TextValue "D38" "1.89"
$ws.Range("E38").Value = "  -2.15%  "
Set-TextValue "D39" "26.81"
$ws.Range("E39").Value = "  -0.62%  "
Set-TextValue "D40" "0.808"
$ws.Range("E40").Value = "  -4.12%  "
Set-TextValue "D41" "4.52"
$ws.Range("E41").Value = "  -2.50%  "
Set-TextValue "D42" "6.42"
$ws.Range("E42").Value = "  -4.08%  "
Set-TextValue "D43" "2.674.59"
$ws.Range("E43").Value = "  +0.09%  "
Set-TextValue "D44" "0.0684"
$ws.Range("E44").Value = "  -2.51%  "
Set-TextValue "D45" "40.73"
$ws.Range("E45").Value = "  -2.78%  "
Set-TextValue "D46" "2.43"
$ws.Range("E46").Value = "  -8.46%  "
Set-TextValue "D47" "24.70"
$ws.Range("E47").Value = "  -4.54%  "
Set-TextValue "D48" "334.68"
$ws.Range("E48").Value = "  -2.66%  "
Set-TextValue "D49" "0.0276"
$ws.Range("E49").Value = "  -4.18%  "
Set-TextValue "D50" "6.37"
$ws.Range("E50").Value = "  +0.29%  "
Set-TextValue "D51" "0.101"
$ws.Range("E51").Value = "  -1.87%  "

$scratch.ClearContents()
$excel.CutCopyMode = 0
Write-Output "Applied 91 cell updates"
